$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.444.71'
$ws.Range('E2').Value = '  +2.28%  '
$ws.Range('D3').Value = '2.005.53'
$ws.Range('E3').Value = '  +5.53%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.19'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.659'
$ws.Range('E6').Value = '  -4.70%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '44.38'
$ws.Range('E8').Value = '  +2.40%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '61.30'
$ws.Range('E9').Value = '  +8.77%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.365'
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0713'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0979'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.51'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '2.298.03'
$ws.Range('E14').Value = '  +5.58%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.808'
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').Value = '2.006.22'
$ws.Range('E16').Value = '  +5.53%  '
$ws.Range('E17').Value = '  -2.92%  '
$ws.Range('D18').Value = '36.312.70'
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.14'
$ws.Range('E19').Value = '  -3.63%  '
$ws.Range('D20').Value = '0.0₃0815'
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.79'
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '236.79'
$ws.Range('E22').Value = '  -3.85%  '
$ws.Range('E23').Value = '  -6.62%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').Value = '  -10.19%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.59'
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.61'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.45'
$ws.Range('E28').Value = '  +5.83%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.91'
$ws.Range('E29').Value = '  -11.72%  '
$ws.Range('E30').Value = '  -5.69%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.30'
$ws.Range('E31').Value = '  +63.88%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.35'
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0582'
$ws.Range('E33').Value = '  -3.60%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.97'
$ws.Range('E36').Value = '  -6.62%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0809'
$ws.Range('E37').Value = '  +9.28%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.12'
$ws.Range('E38').Value = '  +8.33%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.851'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.33'
$ws.Range('E40').Value = '  -11.04%  '
$ws.Range('E41').Value = '  -4.75%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '95.48'
$ws.Range('E42').Value = '  -4.17%  '
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('E44').Value = '  +15.50%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '15.94'
$ws.Range('E45').Value = '  -6.88%  '
$ws.Range('D46').Value = '1.310.85'
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0814'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('D49').Value = '2.192.97'
$ws.Range('E49').Value = '  +5.54%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.18'
$ws.Range('E50').Value = '  -8.06%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.85'
$ws.Range('E51').Value = '  +13.94%  '
